# Apply "Updated symbol list" edits to the cryptos worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while preserving its original
# "text" storage (the source file stores all of these as inline
# strings, so plain numeric-looking assignments must be forced to
# stay textual rather than being reinterpreted as numbers).
function Set-TextValue {
    param($range, [string]$value)
    $range.NumberFormat = "@"
    $range.Value = $value
}

# --- Column D (Price) updates -------------------------------------------------
Set-TextValue $ws.Range("D2")  "275.55"
Set-TextValue $ws.Range("D4")  "6.211"
Set-TextValue $ws.Range("D5")  "0.06180"
Set-TextValue $ws.Range("D6")  "3.575"
Set-TextValue $ws.Range("D8")  "6.529"
Set-TextValue $ws.Range("D10") "0.1645"
Set-TextValue $ws.Range("D11") "0.08261"
Set-TextValue $ws.Range("D12") "0.03432"
Set-TextValue $ws.Range("D14") "0.09137"
Set-TextValue $ws.Range("D15") "3.772"
Set-TextValue $ws.Range("D16") "0.001607"
Set-TextValue $ws.Range("D17") "0.04675"
Set-TextValue $ws.Range("D18") "0.006473"
Set-TextValue $ws.Range("D19") "0.006135"
Set-TextValue $ws.Range("D22") "3.727"
Set-TextValue $ws.Range("D23") "2.322"
Set-TextValue $ws.Range("D24") "0.01385"
Set-TextValue $ws.Range("D25") "0.3278"
Set-TextValue $ws.Range("D26") "0.1230"
Set-TextValue $ws.Range("D40") "0.04742"

# --- Rows 41-43: the three coins rotate position (row 41 now holds what
#     used to be row 42's coin, row 42 holds row 43's coin, and row 43
#     holds row 41's coin), each with its own refreshed price/label. ----------

# Row 41 becomes KickToken
Set-TextValue $ws.Range("B41") "KickToken"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws.Range("D41") "0.007037"
Set-TextValue $ws.Range("E41") "40KickTokenKICK"

# Row 42 becomes BKEXToken
Set-TextValue $ws.Range("B42") "BKEXToken"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws.Range("D42") "0.1104"
Set-TextValue $ws.Range("E42") "41BKEXTokenBKK"

# Row 43 becomes CEJI
Set-TextValue $ws.Range("B43") "CEJI"
Set-TextValue $ws.Range("C43") "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws.Range("D43") "0.005401"
Set-TextValue $ws.Range("E43") "42CEJICEJIBestin24h"

# --- Remaining column D (Price) updates --------------------------------------
Set-TextValue $ws.Range("D44") "0.01149"
Set-TextValue $ws.Range("D48") "0.001387"
